$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.919.54'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '2.442.21'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.172'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.98%  '
$ws.Range('E10').Value = '  -2.19%  '
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.61'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000179'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.64%  '
$ws.Range('D14').Value = '68.764.95'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').Value = '2.886.97'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.31'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').Value = '2.439.56'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.61'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '339.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.95'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.16%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.72'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('D26').Value = '2.567.61'
$ws.Range('E26').Value = '  -1.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.972'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.57%  '
$ws.Range('D29').Value = '0.0₃0822'
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.15'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E32').Value = '  +2.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '430.15'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('E34').Value = '  -1.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '160.07'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.02'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('E39').Value = '  -1.99%  '
$ws.Range('E40').Value = '  +0.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.52'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.37'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.07'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.33'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '130.40'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.482'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.560'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('E50').Value = '  +2.89%  '
